$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4444.3335
$ws.Range("I51").Value = 4999.6665
$ws.Range("J51").Value = 4166.6665
$ws.Range("K51").Value = 4999.6665
$ws.Range("L51").Value = 4166.6665
$ws.Range("M51").Value = -4515.6665
$ws.Range("N51").Value = -5134.6665
$ws.Range("H113").Value = 6150.8667
$ws.Range("I113").Value = 8451.333000000001
$ws.Range("J113").Value = 4617.222
$ws.Range("K113").Value = 8451.333000000001
$ws.Range("L113").Value = 4617.222
$ws.Range("M113").Value = -5197.333000000001
$ws.Range("N113").Value = -11125.222
$ws.Range("H116").Value = 4422
$ws.Range("I116").Value = 3952
$ws.Range("J116").Value = 4735.3335
$ws.Range("K116").Value = 3952
$ws.Range("L116").Value = 4735.3335
$ws.Range("M116").Value = -510
$ws.Range("N116").Value = -11619.3335
$ws.Range("H137").Value = 4259.607
$ws.Range("I137").Value = 5745
$ws.Range("J137").Value = 3764.476
$ws.Range("K137").Value = 17235
$ws.Range("L137").Value = 11293.428
$ws.Range("M137").Value = -14685
$ws.Range("N137").Value = -16393.428
$ws.Range("H138").Value = 38469044
$ws.Range("J138").Value = 13099.2
$ws.Range("L138").Value = 39297.60000000001
$ws.Range("N138").Value = -49577.60000000001
$ws.Range("H141").Value = 7597.2393
$ws.Range("I141").Value = 7343.8647
$ws.Range("J141").Value = 8638.888999999999
$ws.Range("K141").Value = 22031.5941
$ws.Range("L141").Value = 25916.667
$ws.Range("M141").Value = -16851.5941
$ws.Range("N141").Value = -36276.667
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 24999
$ws.Range("J55").Value = 24999
$ws.Range("L55").Value = 24999
$ws.Range("N55").Value = -25629
$ws.Range("H97").Value = 2259.8333
$ws.Range("I97").Value = 2245.5908
$ws.Range("K97").Value = 2245.5908
$ws.Range("M97").Value = -1749.5908
$ws.Range("H102").Value = 5105.25
$ws.Range("I102").Value = 4473.4546
$ws.Range("K102").Value = 4473.4546
$ws.Range("M102").Value = -2851.4546
$ws.Range("H132").Value = 2317.12
$ws.Range("I132").Value = 1706.9
$ws.Range("K132").Value = 5120.700000000001
$ws.Range("M132").Value = -2590.700000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5962
$ws.Range("I99").Value = 5962
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5962
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -4464
$ws.Range("N99").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3523.1738
$ws.Range("I31").Value = 3378.111
$ws.Range("J31").Value = 3616.4285
$ws.Range("K31").Value = 3378.111
$ws.Range("L31").Value = 3616.4285
$ws.Range("M31").Value = -3083.111
$ws.Range("N31").Value = -4206.4285
$ws.Range("H34").Value = 3523.1738
$ws.Range("I34").Value = 3378.111
$ws.Range("J34").Value = 3616.4285
$ws.Range("K34").Value = 3378.111
$ws.Range("L34").Value = 3616.4285
$ws.Range("M34").Value = -3176.111
$ws.Range("N34").Value = -4020.4285
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1200
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1200
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3600
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4520
$ws.Range("H55").Value = 1564.2858
$ws.Range("I55").Value = 1625
$ws.Range("K55").Value = 4875
$ws.Range("M55").Value = -4698
$ws.Range("H86").Value = 1220.6364
$ws.Range("I86").Value = 274.6
$ws.Range("J86").Value = 2009
$ws.Range("K86").Value = 823.8000000000001
$ws.Range("L86").Value = 6027
$ws.Range("M86").Value = 362.1999999999999
$ws.Range("N86").Value = -8399
$ws.Range("H89").Value = 1220.6364
$ws.Range("I89").Value = 274.6
$ws.Range("J89").Value = 2009
$ws.Range("K89").Value = 2471.4
$ws.Range("L89").Value = 18081
$ws.Range("M89").Value = 3456.6
$ws.Range("N89").Value = -29937
$ws.Range("H140").Value = 1878.8
$ws.Range("I140").Value = 1098.625
$ws.Range("K140").Value = 3295.875
$ws.Range("M140").Value = 1884.125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10303.875
$ws.Range("I70").Value = 12821.637
$ws.Range("J70").Value = 4764.8
$ws.Range("K70").Value = 12821.637
$ws.Range("L70").Value = 4764.8
$ws.Range("M70").Value = -12551.637
$ws.Range("N70").Value = -5304.8
$ws.Range("H73").Value = 10303.875
$ws.Range("I73").Value = 12821.637
$ws.Range("J73").Value = 4764.8
$ws.Range("K73").Value = 12821.637
$ws.Range("L73").Value = 4764.8
$ws.Range("M73").Value = -11885.637
$ws.Range("N73").Value = -6636.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3483.5
$ws.Range("I7").Value = 2778.375
$ws.Range("J7").Value = 4423.6665
$ws.Range("K7").Value = 2778.375
$ws.Range("L7").Value = 4423.6665
$ws.Range("M7").Value = -2666.375
$ws.Range("N7").Value = -4647.6665
$ws.Range("H61").Value = 6148.375
$ws.Range("J61").Value = 2945.6667
$ws.Range("L61").Value = 2945.6667
$ws.Range("N61").Value = -3349.6667
$ws.Range("H93").Value = 2375.8462
$ws.Range("I93").Value = 2681.0908
$ws.Range("J93").Value = 697
$ws.Range("K93").Value = 2681.0908
$ws.Range("L93").Value = 697
$ws.Range("M93").Value = -1433.0908
$ws.Range("N93").Value = -3193
$ws.Range("H113").Value = 6148.375
$ws.Range("J113").Value = 2945.6667
$ws.Range("L113").Value = 2945.6667
$ws.Range("N113").Value = -7285.6667
$ws.Range("H126").Value = 3483.5
$ws.Range("I126").Value = 2778.375
$ws.Range("J126").Value = 4423.6665
$ws.Range("K126").Value = 8335.125
$ws.Range("L126").Value = 13270.9995
$ws.Range("M126").Value = -5865.125
$ws.Range("N126").Value = -18210.9995
$ws.Range("H132").Value = 5340
$ws.Range("I132").Value = 4900
$ws.Range("K132").Value = 14700
$ws.Range("M132").Value = -12170
$ws.Range("H136").Value = 8189185
$ws.Range("I136").Value = 10590540
$ws.Range("J136").Value = 24577.4
$ws.Range("K136").Value = 31771620
$ws.Range("L136").Value = 73732.20000000001
$ws.Range("M136").Value = -31769070
$ws.Range("N136").Value = -78832.20000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4893.3823
$ws.Range("I132").Value = 4732.684
$ws.Range("K132").Value = 14198.052
$ws.Range("M132").Value = -11668.052
